$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: name and week number
$ws.Range("C1").Value = "Richard Dobson"
$ws.Range("E1").Value = 12

# Row 3: first task entry (also drop the stray wrap-text formatting so its
# style matches the rest of the rows, as in the target workbook)
$ws.Range("A3").WrapText = $False
$ws.Range("A3").Value = "Project Build"
$ws.Range("B3").Value = "Prepare Client Presentation"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 10

# Row 4: second task entry
$ws.Range("A4").Value = "Project Build"
$ws.Range("B4").Value = "Final Cehcks and tweaks"
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 10

# Update selection to match the saved view state
[void]$ws.Range("D4").Select()
